$d = $word.ActiveDocument

# ======================================================================
# Change 1: paragraph ending in "{Assinatura Militar Arrolado}"
#   - indentation: left 2832->2124 twips, firstLine 708->0 twips
#     (141.6pt -> 106.2pt ; 35.4pt -> 0pt)
#   - remove the leading " " run and the following <w:tab/> run
# ======================================================================
$find1 = $d.Content.Find
$find1.ClearFormatting()
$found1 = $find1.Execute("Assinatura Militar Arrolado", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not locate paragraph containing 'Assinatura Militar Arrolado'"
}
$para1 = $find1.Parent.Paragraphs(1)

# Update paragraph indentation (2124 twips = 106.2 pt, 0 twips = 0 pt)
$para1.Format.LeftIndent = 106.2
$para1.Format.FirstLineIndent = 0

# The paragraph currently starts with a run containing a single space,
# then a run containing just a tab character, then the real text run.
# Delete those first two characters (space + tab) to remove both runs.
$para1Range = $para1.Range
$leadRange = $d.Range($para1Range.Start, $para1Range.Start + 2)
if ($leadRange.Text -eq " `t") {
    $leadRange.Delete()
}

# ======================================================================
# Change 2: paragraph containing "PATD N..." / the PATD number
#   - pPr/rPr: bold on, language en-US
#   - first run: fewer leading spaces (75 -> 70), language en-US
#   - space run: language en-US
#   - final run "0947/BAGL-GSDGL/18072025" is replaced by the
#     placeholder sequence "{N PATD}/BAGL-GSDGL/{DataPatd}" spread
#     across 7 runs (first one at sz 24, the rest at sz 22)
# ======================================================================
$find2 = $d.Content.Find
$find2.ClearFormatting()
$found2 = $find2.Execute("0947/BAGL-GSDGL/18072025", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not locate the PATD number run"
}
$numStart = $find2.Parent.Start
$numEnd = $find2.Parent.End
$numRange = $d.Range($numStart, $numEnd)
$para2 = $numRange.Paragraphs(1)
$para2Range = $para2.Range

# --- Bold the whole paragraph (propagates w:b to the mark's rPr too) ---
$para2Range.Font.Bold = $true
$para2Range.Font.BoldBi = $true
$para2Range.Font.LanguageID = "en-US"

# --- Shrink the run of leading spaces before "PATD N" by 5 spaces ---
$findSpaces = $d.Content.Find
$findSpaces.ClearFormatting()
$findSpaces.Replacement.ClearFormatting()
$spacesFound = $findSpaces.Execute("                                                                           PATD N", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "                                                                      PATD N", 2)
if (-not $spacesFound) {
    throw "Could not shrink the leading spaces before PATD N"
}

# --- Replace the static number run with the templated placeholder runs ---
$find3 = $d.Content.Find
$find3.ClearFormatting()
$found3 = $find3.Execute("0947/BAGL-GSDGL/18072025", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) {
    throw "Could not re-locate the PATD number run for replacement"
}
$targetStart = $find3.Parent.Start
$targetEnd = $find3.Parent.End
$targetRange = $d.Range($targetStart, $targetEnd)

$newRunsXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="1"/><w:bCs w:val="1"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:noProof w:val="0"/><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>{N PATD}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="1"/><w:bCs w:val="1"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:noProof w:val="0"/><w:color w:val="FF0000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>/BAGL</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="1"/><w:bCs w:val="1"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:noProof w:val="0"/><w:color w:val="FF0000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>-</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="1"/><w:bCs w:val="1"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:noProof w:val="0"/><w:color w:val="FF0000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>GSDGL/</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="1"/><w:bCs w:val="1"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:noProof w:val="0"/><w:color w:val="FF0000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="1"/><w:bCs w:val="1"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:noProof w:val="0"/><w:color w:val="FF0000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>DataPatd</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="1"/><w:bCs w:val="1"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:noProof w:val="0"/><w:color w:val="FF0000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($newRunsXml)
